$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 636.0714
$ws.Range("I33").Value = 536.9091
$ws.Range("J33").Value = 999.6667
$ws.Range("K33").Value = 536.9091
$ws.Range("L33").Value = 999.6667
$ws.Range("M33").Value = -307.9091
$ws.Range("N33").Value = -1457.6667
$ws.Range("H93").Value = 32499.2
$ws.Range("J93").Value = 32499.2
$ws.Range("L93").Value = 32499.2
$ws.Range("N93").Value = -37491.2
$ws.Range("H100").Value = 3340
$ws.Range("I100").Value = 3340
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 3340
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -2799
$ws.Range("N100").ClearContents()
$ws.Range("H125").Value = 1300.8
$ws.Range("I125").Value = 1212.5454
$ws.Range("J125").Value = 1408.6666
$ws.Range("K125").Value = 10912.9086
$ws.Range("L125").Value = 12677.9994
$ws.Range("M125").Value = -8452.908599999999
$ws.Range("N125").Value = -17597.9994
$ws.Range("H129").Value = 1751.9524
$ws.Range("I129").Value = 1832.1
$ws.Range("J129").Value = 1679.091
$ws.Range("K129").Value = 5496.299999999999
$ws.Range("L129").Value = 5037.272999999999
$ws.Range("M129").Value = -496.2999999999993
$ws.Range("N129").Value = -15037.273

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H53").Value = 12679.667
$ws.Range("I53").Value = 8519.5
$ws.Range("J53").Value = 21000
$ws.Range("K53").Value = 8519.5
$ws.Range("L53").Value = 21000
$ws.Range("M53").Value = -7837.5
$ws.Range("N53").Value = -22364
$ws.Range("H61").Value = 2703.4255
$ws.Range("I61").Value = 1546.1034
$ws.Range("K61").Value = 1546.1034
$ws.Range("M61").Value = -1334.1034
$ws.Range("H101").Value = 46596
$ws.Range("J101").Value = 46596
$ws.Range("L101").Value = 46596
$ws.Range("N101").Value = -53086
$ws.Range("H103").Value = 39354
$ws.Range("J103").Value = 39354
$ws.Range("L103").Value = 39354
$ws.Range("N103").Value = -41698
$ws.Range("H106").Value = 47681
$ws.Range("J106").Value = 47681
$ws.Range("L106").Value = 47681
$ws.Range("N106").Value = -50205
$ws.Range("H117").Value = 47311.75
$ws.Range("J117").Value = 47311.75
$ws.Range("L117").Value = 47311.75
$ws.Range("N117").Value = -56489.75
$ws.Range("H136").Value = 2703.4255
$ws.Range("I136").Value = 1546.1034
$ws.Range("K136").Value = 4638.3102
$ws.Range("M136").Value = -2088.3102

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1906.3334
$ws.Range("I105").Value = 1239.1666
$ws.Range("J105").Value = 3240.6667
$ws.Range("K105").Value = 1239.1666
$ws.Range("L105").Value = 3240.6667
$ws.Range("M105").Value = 507.8334
$ws.Range("N105").Value = -6734.6667

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1654.5883
$ws.Range("I58").Value = 1409.4642
$ws.Range("J58").Value = 2798.5
$ws.Range("K58").Value = 1409.4642
$ws.Range("L58").Value = 2798.5
$ws.Range("M58").Value = -1206.4642
$ws.Range("N58").Value = -3204.5
$ws.Range("H125").Value = 30563
$ws.Range("J125").Value = 30563
$ws.Range("L125").Value = 30563
$ws.Range("N125").Value = -35483
$ws.Range("H131").Value = 38326
$ws.Range("J131").Value = 38326
$ws.Range("L131").Value = 38326
$ws.Range("N131").Value = -48406
$ws.Range("H136").Value = 1654.5883
$ws.Range("I136").Value = 1409.4642
$ws.Range("J136").Value = 2798.5
$ws.Range("K136").Value = 4228.392599999999
$ws.Range("L136").Value = 8395.5
$ws.Range("M136").Value = -1678.392599999999
$ws.Range("N136").Value = -13495.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5417.0835
$ws.Range("I80").Value = 5775.625
$ws.Range("J80").Value = 4700
$ws.Range("K80").Value = 5775.625
$ws.Range("L80").Value = 4700
$ws.Range("M80").Value = -4777.625
$ws.Range("N80").Value = -6696
$ws.Range("H83").Value = 5417.0835
$ws.Range("I83").Value = 5775.625
$ws.Range("J83").Value = 4700
$ws.Range("K83").Value = 28878.125
$ws.Range("L83").Value = 23500
$ws.Range("M83").Value = -23886.125
$ws.Range("N83").Value = -33484
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H120").Value = 22712.75
$ws.Range("J120").Value = 22712.75
$ws.Range("L120").Value = 22712.75
$ws.Range("N120").Value = -32388.75
$ws.Range("H127").Value = 41996
$ws.Range("J127").Value = 41996
$ws.Range("L127").Value = 41996
$ws.Range("N127").Value = -51916
$ws.Range("H134").Value = 25977.777
$ws.Range("J134").Value = 25977.777
$ws.Range("L134").Value = 77933.33099999999
$ws.Range("N134").Value = -83003.33099999999
$ws.Range("H135").Value = 46318.4
$ws.Range("J135").Value = 46318.4
$ws.Range("L135").Value = 46318.4
$ws.Range("N135").Value = -56458.4

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3302.1052
$ws.Range("I136").Value = 2718.923
$ws.Range("J136").Value = 4565.6665
$ws.Range("K136").Value = 8156.768999999999
$ws.Range("L136").Value = 13696.9995
$ws.Range("M136").Value = -5606.768999999999
$ws.Range("N136").Value = -18796.9995

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 28906.5
$ws.Range("J27").Value = 28906.5
$ws.Range("L27").Value = 28906.5
$ws.Range("N27").Value = -29044.5
$ws.Range("H80").Value = 35483.8
$ws.Range("J80").Value = 35483.8
$ws.Range("L80").Value = 35483.8
$ws.Range("N80").Value = -37479.8
$ws.Range("H83").Value = 35483.8
$ws.Range("J83").Value = 35483.8
$ws.Range("L83").Value = 106451.4
$ws.Range("N83").Value = -116435.4
$ws.Range("H108").Value = 38957.332
$ws.Range("J108").Value = 38957.332
$ws.Range("L108").Value = 38957.332
$ws.Range("N108").Value = -46637.332
$ws.Range("H109").Value = 34614.75
$ws.Range("J109").Value = 34614.75
$ws.Range("L109").Value = 34614.75
$ws.Range("N109").Value = -37388.75
$ws.Range("H115").Value = 37369
$ws.Range("J115").Value = 37369
$ws.Range("L115").Value = 37369
$ws.Range("N115").Value = -40503
$ws.Range("H118").Value = 34586.668
$ws.Range("J118").Value = 44380
$ws.Range("L118").Value = 44380
$ws.Range("N118").Value = -47694
$ws.Range("H129").Value = 29028
$ws.Range("J129").Value = 29028
$ws.Range("L129").Value = 29028
$ws.Range("N129").Value = -39028
$ws.Range("H131").Value = 43575.168
$ws.Range("J131").Value = 43575.168
$ws.Range("L131").Value = 43575.168
$ws.Range("N131").Value = -53655.168
$ws.Range("H136").Value = 22220.041
$ws.Range("I136").Value = 67620.47
$ws.Range("J136").Value = 2190.4412
$ws.Range("K136").Value = 202861.41
$ws.Range("L136").Value = 6571.323600000001
$ws.Range("M136").Value = -200311.41
$ws.Range("N136").Value = -11671.3236
